$wb = $excel.ActiveWorkbook

# ---------- Sheet "LP1912": append rows 948-964, update header ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$s1rows = @(
    @(948, '12:56:02', '12:58', '16_SANTA ANA', 2, 'LP1912', '31/12/2025'),
    @(949, '12:56:02', '12:59', '10_OLMOS', 3, 'LP1912', '31/12/2025'),
    @(950, '12:56:02', '13:01', '215C_EL PATO', 5, 'LP1912', '31/12/2025'),
    @(951, '12:56:02', '13:04', '23_HERNANDEZ', 8, 'LP1912', '31/12/2025'),
    @(952, '12:56:02', '13:07', '14_ABASTO', 11, 'LP1912', '31/12/2025'),
    @(953, '12:56:02', '13:11', '16_SANTA ANA', 15, 'LP1912', '31/12/2025'),
    @(954, '12:56:02', '13:19', '11_ETCHEVERRY', 23, 'LP1912', '31/12/2025'),
    @(955, '12:56:02', '13:21', '16_SANTA ANA', 25, 'LP1912', '31/12/2025'),
    @(956, '12:56:02', '13:21', '17_ROMERO', 25, 'LP1912', '31/12/2025'),
    @(957, '12:56:02', '13:30', '10_OLMOS', 34, 'LP1912', '31/12/2025'),
    @(958, '12:56:02', '13:31', '16_P MOR-SANTA ANA', 35, 'LP1912', '31/12/2025'),
    @(959, '12:56:02', '13:34', '23_HERNANDEZ', 38, 'LP1912', '31/12/2025'),
    @(960, '12:56:02', '13:51', '15_ABASTO', 55, 'LP1912', '31/12/2025'),
    @(961, '12:56:02', '14:01', '17_ROMERO', 65, 'LP1912', '31/12/2025'),
    @(962, '12:56:02', '14:04', '23_HERNANDEZ', 68, 'LP1912', '31/12/2025'),
    @(963, '12:56:02', '14:11', '15_ABASTO', 75, 'LP1912', '31/12/2025'),
    @(964, '12:56:02', '14:25', '11_ETCHEVERRY', 89, 'LP1912', '31/12/2025')
)
foreach ($r in $s1rows) {
    $row = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $ws1.Cells.Item($row, 6).Value = $r[5]
    $ws1.Cells.Item($row, 7).Value = $r[6]
}
$ws1.Range("A2").Value = "Última actualización: 31/12/2025 12:56:12"
$ws1.Range("A3").Value = "Total filas: 963"

# ---------- Sheet "LP1912-215": append row 74, update header ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(74, 2).Value = '31/12/2025'
$ws2.Cells.Item(74, 3).Value = '12:56:02'
$ws2.Cells.Item(74, 4).Value = '13:01'
$ws2.Cells.Item(74, 5).Value = '215C_EL PATO'
$ws2.Cells.Item(74, 6).Value = 5
$ws2.Cells.Item(74, 7).Value = 'LP1912'
$ws2.Range("A2").Value = "Última actualización: 31/12/2025 12:56:12"
$ws2.Range("A3").Value = "Total filas: 73"

# ---------- Sheet "6203-6173": append rows 115-118, update header ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$s3rows = @(
    @(115, '31/12/2025', '12:56:12', '13:09', '215B_LP-P MOR-1 Y 57', 13, 'L6173'),
    @(116, '31/12/2025', '12:56:12', '13:14', '215A_LA PLATA', 18, 'L6173'),
    @(117, '31/12/2025', '12:56:07', '13:54', '215C_LA PLATA', 58, 'L6203'),
    @(118, '31/12/2025', '12:56:07', '14:34', '215C_LA PLATA', 98, 'L6203')
)
foreach ($r in $s3rows) {
    $row = $r[0]
    $ws3.Cells.Item($row, 2).Value = $r[1]
    $ws3.Cells.Item($row, 3).Value = $r[2]
    $ws3.Cells.Item($row, 4).Value = $r[3]
    $ws3.Cells.Item($row, 5).Value = $r[4]
    $ws3.Cells.Item($row, 6).Value = $r[5]
    $ws3.Cells.Item($row, 7).Value = $r[6]
}
$ws3.Range("A2").Value = "Última actualización: 31/12/2025 12:56:12"
$ws3.Range("A3").Value = "Total filas: 117"
